# Updated README & PowerPoint
#
# Slide 4 ("Take Aways"), "Content Placeholder 2" shape: the "For Next
# Time" list gets expanded with extra bullet points around the existing
# "technicality" line (new lines before it, a blank separator line, and
# new lines after it), while the original "technicality" run and the
# trailing endParaRPr on the box stay intact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 2 is currently "technicality" (paragraph 1 is "For Next Time").
# Insert the two new bullets plus a blank line right before it - this
# pushes "technicality" down to paragraph 5, and keeps the shape's
# trailing <a:endParaRPr lang="de-CH"/> attached to it (it travels with
# the original paragraph mark).
$beforePara = $tr.Paragraphs(2, 1)
$beforePara.InsertBefore("- Improve Time management`r- better work distribution`r`r")

# Re-fetch the range/paragraph after the structural change and append the
# two closing bullets after "technicality" (now paragraph 5) - this moves
# the trailing endParaRPr to the very last new paragraph, matching the
# target structure.
$fresh = $shape.TextFrame.TextRange
$techPara = $fresh.Paragraphs(5, 1)
$techPara.InsertAfter("`r- Work on " + [char]0x2018 + "dev" + [char]0x2019 + " branch`r- know your skills")
